# "updated canine keywords; added samples, files obj"
# Update the Case ID / Diagnosis / Stage of Disease values on the
# CypherOutput results sheet (the active sheet, row 2 = the data row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "NCATS-COP01-CCB010072"
$ws.Range("E2").Value = "Osteosarcoma"
$ws.Range("F2").Value = "Unknown"
